# CreateEntity_OOFS_MultiSuccessCase_Test.xlsx — "create entity in progress"
#
# The fixture is a log of an automated "Create Entity" workflow run. Re-running
# the workflow moved on to the next day's pass (2024-04-04 -> 2024-04-05) and
# refreshed the per-row "started at" / "confirmed at" timestamp columns on the
# first two result sheets, and left the UI focused on Sheet2 with a different
# active cell on each sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- refresh the timestamp columns on Sheet1 & Sheet2 ------------------------
# These columns hold plain text that merely looks like dates/times (e.g.
# "2024-04-05", "CT: Fri, Apr 05, 2024 at 7:45 PM"). Assigning such a string
# straight to .Value gets auto-coerced into a real date serial by the COM
# layer, which would also swap in a date number format — neither of which
# matches the source data (plain text, General format). Stage the literal
# text in a scratch cell that's forced to Text format, then copy only the
# value across so the destination keeps its original (General) style.
$scratch = $ws1.Range("ZZ100")
$scratch.NumberFormat = "@"

function Set-LiteralText($range, $text) {
    $scratch.Value = $text
    $scratch.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

Set-LiteralText $ws1.Range("O2")  "2024-04-05"
Set-LiteralText $ws1.Range("Q2")  "2024-04-05 07:43:05 PM"
Set-LiteralText $ws1.Range("AD2") "2024-04-05"
Set-LiteralText $ws1.Range("BB2") "CT: Fri, Apr 05, 2024 at 7:45 PM"

Set-LiteralText $ws2.Range("O2")  "2024-04-05"
Set-LiteralText $ws2.Range("Q2")  "2024-04-05 07:52:58 PM"
Set-LiteralText $ws2.Range("AD2") "2024-04-05"
Set-LiteralText $ws2.Range("BB2") "CT: Fri, Apr 05, 2024 at 8:02 PM"

$scratch.Clear()
$excel.CutCopyMode = $false

# --- move the active selection / active sheet --------------------------------
# Sheet1's selection moves to AJ2 (and it stops being the active tab); Sheet2
# becomes the active tab with its selection moved to E2.
$ws1.Range("AJ2").Select() | Out-Null
$ws2.Range("E2").Select() | Out-Null
$ws2.Activate() | Out-Null
